$d = $word.ActiveDocument

$pairs = @(
    @("882×2=", "766×4="),
    @("691×5=", "811×7="),
    @("411×5=", "801×9="),
    @("496×2=", "345×9="),
    @("542×6=", "572×2="),
    @("869×6=", "904×7="),
    @("595×3=", "332×5="),
    @("878×6=", "575×7="),
    @("421×5=", "914×4="),
    @("663×9=", "704×9="),
    @("465×4=", "420×3="),
    @("637×5=", "866×4="),
    @("347×4=", "624×5="),
    @("274×6=", "251×2="),
    @("350×5=", "352×3="),
    @("435×3=", "122×3="),
    @("390×4=", "432×5="),
    @("200×4=", "883×2="),
    @("949×3=", "765×9="),
    @("633×4=", "308×5="),
    @("745×8=", "378×4="),
    @("764×7=", "684×6="),
    @("218×6=", "359×3="),
    @("519×6=", "899×6="),
    @("569×8=", "230×9=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
